# A2osX Issue List update
# - Mark several resolved Issues rows as Closed ("C")
# - Add 3 new Issue rows (IFNUMTEST related)
# - Add / complete several Suggestions rows (S0048-S0055)

$wb = $excel.ActiveWorkbook

$wsIssues = $wb.Worksheets.Item("Issues")
$wsSugg   = $wb.Worksheets.Item("Suggestions")

# ---------------------------------------------------------------
# Suggestions sheet: finish off S0048, add S0049-S0055
# (edited in this order so new shared strings line up the same
#  way the author originally typed them)
# ---------------------------------------------------------------

$wsSugg.Range("D49").Value = 'Editing comand line, Ctrl-X deletes rest of line, Ctrl-C clear line but keep in history, hist is 256 bytes'
$wsSugg.Range("C49").Value = 'Doc'
$wsSugg.Rows(49).RowHeight = 30

$wsSugg.Range("B50").Value = 'S0049'
$wsSugg.Range("C50").Value = 'SHELL'
$wsSugg.Range("D50").Value = 'Add the ability to parse substrings ${VAR:12:22}.  This will likely also require a LEN($VAR) command.'
$wsSugg.Rows(50).RowHeight = 30

$wsSugg.Range("B51").Value = 'S0050'
$wsSugg.Range("C51").Value = 'ECHO'

$wsSugg.Range("B52").Value = 'S0051'
$wsSugg.Range("B53").Value = 'S0052'
$wsSugg.Range("B54").Value = 'S0053'
$wsSugg.Range("B55").Value = 'S0054'
$wsSugg.Range("B56").Value = 'S0055'

$wsSugg.Range("D51").Value = 'Add a way to ECHO "something" with no CR/LF, so that the next ECHO "else" will put on screen "somethingelse".  There are many uses like Displaying "Test IF LT: " and then after test run "Pass or Fail" on same line.  Or if processing large loop, can ECHO "."; and have dots appear for each loop run.'
$wsSugg.Rows(51).RowHeight = 90

# ---------------------------------------------------------------
# Issues sheet: mark closed issues
# ---------------------------------------------------------------

$wsIssues.Range("A55").Value = 'C'
$wsIssues.Range("A56").Value = 'C'
$wsIssues.Range("A58").Value = 'C'
$wsIssues.Range("A63").Value = 'C'
$wsIssues.Range("A65").Value = 'C'
$wsIssues.Range("A66").Value = 'C'
$wsIssues.Range("A70").Value = 'C'

# ---------------------------------------------------------------
# Issues sheet: new rows for IFNUMTEST related issues
# ---------------------------------------------------------------

$wsIssues.Range("B72").Value = 164
$wsIssues.Range("C72").Value = 'SHELL'
$wsIssues.Range("D72").Value = 1079
$wsIssues.Range("E72").Value = 'If you are running a long script (IFNUMTEST) and press Ctrl-S to pause output you get [$80} Unknown Error in the script, and you have to press return to continue the script.'
$wsIssues.Rows(72).RowHeight = 60

$wsIssues.Range("B73").Value = 165
$wsIssues.Range("C73").Value = 'EDIT'
$wsIssues.Range("D73").Value = 1079
$wsIssues.Range("E73").Value = 'EDIT longfile.  Down arrow 5 or 6 times, hit Ctrl-P to page down, the top 5-6 lines get blanked and cursor goes bottom.  If you Ctrl-O/P at top/bttm screen its ok'
$wsIssues.Rows(73).RowHeight = 60

$wsIssues.Range("B74").Value = 166
$wsIssues.Range("C74").Value = 'IF Numerics'
$wsIssues.Range("D74").Value = 1079
$wsIssues.Range("E74").Value = 'Several Tests failing, see new IFNUMTEST tests 4 8 9 12 and 14'
$wsIssues.Rows(74).RowHeight = 30
